$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: "Rectangle 108" background rectangle - reposition/resize
# a:off x="1" y="0" -> x="1" y="1062032"   (EMU, 12700 EMU per point)
# a:ext cx="12192000" cy="6857999" -> cx="12192000" cy="5194929"
$rect = $s.Shapes.Item(1)
$rect.Top = 83.62457
$rect.Height = 409.04953

# Shape 4: "TextBox 107" - text "Decent" -> "Before"
$s.Shapes.Item(4).TextFrame.TextRange.Text = "Before"

# Shape 5: "TextBox 1" - text "Great" -> "After"
$s.Shapes.Item(5).TextFrame.TextRange.Text = "After"
